# Add a new data row (row 25) for the AWS "i3.16xlarge" storage-optimized
# instance type, which has no Azure equivalent (marked "NA" in the
# Azure-related columns). Cells are written in the same left-to-right
# column order that the new shared strings appear in the target workbook
# (i3.16xlarge, then $4.9920/hour, then NA) so the shared-string table
# indices line up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "Storage Optimized"
$ws.Range("C25").Value = "i3.16xlarge"
$ws.Range("G25").Value = "$4.9920/hour"
$ws.Range("B25").Value = "NA"
$ws.Range("D25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("E25").Value = "64 vCPU | 512 GB RAM"

# Re-apply explicit ("best fit") column widths now that the table has one
# more row of data, matching the widths Excel computed when it re-saved
# the workbook.
$ws.Columns.Item(1).ColumnWidth = 16.417666666666666
$ws.Columns.Item(2).ColumnWidth = 17.417666666666666
$ws.Columns.Item(3).ColumnWidth = 16.584333333333333
$ws.Columns.Item(4).ColumnWidth = 25.084333333333333
$ws.Columns.Item(5).ColumnWidth = 24.251
$ws.Columns.Item(6).ColumnWidth = 11.251
$ws.Columns.Item(7).ColumnWidth = 11.251

# The saved workbook's active selection moved to L17.
$ws.Range("L17").Select() | Out-Null
